# Automatic update of files.
# Bump the "Förändrad" (Changed) date column (C2:C34) by one day:
# from serial date 45639 (2024-12-13) to 45640 (2024-12-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 34; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45639) {
        $cell.Value2 = 45640
    }
}
